# Scenario_Control.xlsx - update Execution Status flags and selection
# Commit: "Removed scrollintoview and scrollby functions to another
#          function library WebObjectAction.qfl"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenario")

# AddToCart / Purchase scenario block (rows 2-3): YES -> No
$ws.Range("F2").Value = "No"
$ws.Range("F3").Value = "No"

# Login / Account scenario block (rows 6-7): Yes -> No
$ws.Range("F6").Value = "No"
$ws.Range("F7").Value = "No"

# Contact Us scenario block (rows 25-26): No -> Yes
$ws.Range("F25").Value = "Yes"
$ws.Range("F26").Value = "Yes"

# Update the active selection to match the author's final cursor position
$ws.Range("J20").Select()
